$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column DL, shifting DL:MN one column to the
# right (DL:MN -> DM:MO). This adds the new "DemonstrationProjectIdentifier"
# field right before the existing "Note" column.
$ws.Range("DL1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("DL1").Value = "DemonstrationProjectIdentifier"

# Update the record identifiers on the two data rows.
$ws.Range("A2").Value = "6901488a7e79911955eafdd4"
$ws.Range("A3").Value = "6901488a7e79911955eafdd4"
